# Applies the commit: "feat: add 2022-Q1 data"
# - The old "总计" sheet (6th tab) is renamed "2022-Q1" and repopulated
#   with the per-fund holdings table for that quarter.
# - A brand-new "总计" (summary) sheet is appended at the end with the
#   updated totals table, including the new 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: the old "总计" sheet (6th tab) becomes the new "2022-Q1" sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$q1Cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $q1Headers.Length; $i++) {
    $cell = $q1.Range($q1Cols[$i] + "1")
    $cell.Value = $q1Headers[$i]
    $cell.Style = $q1.Range("B1").Style
}

# Each row: @(index, code, name, fundScale, stockPosition, positionPct, heldValueText, heldValueIsZero, positionRank)
$funds = @(
    @(0, "'720001", "'财通价值动量混合", "'38.89", "'79.44", "'7.76", "'3.0179", 0, 3),
    @(1, "'270006", "'广发策略优选混合", "'40.17", "'93.86", "'5.47", "'2.1973", 0, 6),
    @(2, "'001480", "'财通成长优选混合", "'16.12", "'86.02", "'8.16", "'1.3154", 0, 2),
    @(3, "'008188", "'前海开源稳健增长三年持有期混合", "'25.26", "'67.84", "'4.33", "'1.0938", 0, 7),
    @(4, "'159865", "'国泰中证畜牧养殖ETF", "'24.32", "'99.29", "'3.87", "'0.9412", 0, 8),
    @(5, "'501085", "'财通科创主题 3 年封闭运作灵活配置混合", "'14.52", "'89.20", "'4.82", "'0.6999", 0, 9),
    @(6, "'009970", "'财通内需增长12个月定期开放混合", "'21.90", "'64.61", "'3.10", "'0.6789", 0, 2),
    @(7, "'164403", "'前海开源沪港深农业主题精选灵活配置混合（LOF）", "'6.21", "'92.31", "'8.44", "'0.5241", 0, 2),
    @(8, "'003857", "'前海开源周期优选灵活配置混合A", "'5.52", "'90.41", "'8.73", "'0.4819", 0, 2),
    @(9, "'001178", "'前海开源再融资主题精选股票", "'8.08", "'84.20", "'4.37", "'0.3531", 0, 8),
    @(10, "'009062", "'财通智慧成长混合A", "'4.27", "'86.88", "'7.96", "'0.3399", 0, 4),
    @(11, "'501046", "'财通多策略福鑫定期开放灵活配置混合", "'3.26", "'86.41", "'7.76", "'0.2530", 0, 4),
    @(12, "'460007", "'华泰柏瑞行业领先混合", "'3.45", "'94.37", "'7.23", "'0.2494", 0, 5),
    @(13, "'001398", "'华泰柏瑞健康生活灵活配置混合", "'3.09", "'94.36", "'7.10", "'0.2194", 0, 5),
    @(14, "'159867", "'鹏华中证畜牧养殖ETF", "'5.61", "'97.87", "'3.80", "'0.2132", 0, 8),
    @(15, "'003858", "'前海开源周期优选灵活配置混合C", "'2.32", "'90.41", "'8.73", "'0.2025", 0, 2),
    @(16, "'011111", "'华泰柏瑞行业严选混合型证券投资基金A", "'2.36", "'94.67", "'7.81", "'0.1843", 0, 3),
    @(17, "'006775", "'前海开源优质成长混合", "'3.68", "'68.45", "'4.40", "'0.1619", 0, 7),
    @(18, "'011287", "'前海开源聚慧三年持有期混合", "'3.25", "'68.77", "'4.37", "'0.1420", 0, 6),
    @(19, "'009063", "'财通智慧成长混合C", "'1.67", "'86.88", "'7.96", "'0.1329", 0, 4),
    @(20, "'501015", "'财通多策略升级混合（LOF）", "'2.63", "'94.52", "'4.46", "'0.1173", 0, 10),
    @(21, "'010418", "'财通景气行业一年封闭运作混合", "'3.26", "'86.25", "'2.94", "'0.0958", 0, 7),
    @(22, "'006216", "'前海开源价值成长灵活配置混合A", "'1.97", "'67.45", "'4.40", "'0.0867", 0, 6),
    @(23, "'590006", "'中邮中小盘灵活配置混合", "'2.39", "'78.54", "'2.93", "'0.0700", 0, 5),
    @(24, "'002407", "'前海开源恒远灵活配置混合", "'1.27", "'60.12", "'4.39", "'0.0558", 0, 6),
    @(25, "'516760", "'平安中证畜牧养殖ETF", "'1.45", "'97.82", "'3.80", "'0.0551", 0, 8),
    @(26, "'005851", "'财通新视野灵活配置混合A", "'0.87", "'85.27", "'5.33", "'0.0464", 0, 8),
    @(27, "'970032", "'东海证券海睿进取灵活配置混合型集合资产管理计划A", "'1.57", "'78.62", "'2.88", "'0.0452", 0, 9),
    @(28, "'516670", "'招商中证畜牧养殖ETF", "'1.06", "'98.73", "'3.83", "'0.0406", 0, 8),
    @(29, "'001601", "'鑫元鑫新收益灵活配置混合A", "'0.72", "'86.90", "'4.87", "'0.0351", 0, 4),
    @(30, "'970033", "'东海证券海睿进取灵活配置混合型集合资产管理计划B", "'1.18", "'78.62", "'2.88", "'0.0340", 0, 9),
    @(31, "'006217", "'前海开源价值成长灵活配置混合C", "'0.73", "'67.45", "'4.40", "'0.0321", 0, 6),
    @(32, "'161715", "'招商中证大宗商品股票指数（LOF）", "'2.24", "'95.07", "'1.40", "'0.0314", 0, 5),
    @(33, "'005959", "'财通新视野灵活配置混合C", "'0.40", "'85.27", "'5.33", "'0.0213", 0, 8),
    @(34, "'011112", "'华泰柏瑞行业严选混合型证券投资基金C", "'0.26", "'94.67", "'7.81", "'0.0203", 0, 3),
    @(35, "'003513", "'中邮消费升级灵活配置混合", "'0.41", "'54.24", "'4.54", "'0.0186", 0, 6),
    @(36, "'010951", "'中金成长精选混合A", "'0.56", "'88.74", "'2.91", "'0.0163", 0, 10),
    @(37, "'970046", "'东海证券海睿健行灵活配置混合A", "'0.57", "'77.45", "'2.85", "'0.0162", 0, 8),
    @(38, "'000417", "'国联安新精选灵活配置混合", "'0.52", "'28.16", "'2.60", "'0.0135", 0, 1),
    @(39, "'003981", "'中银证券瑞益灵活配置混合C", "'0.21", "'89.21", "'5.18", "'0.0109", 0, 3),
    @(40, "'570006", "'诺德中小盘混合", "'0.20", "'90.18", "'3.57", "'0.0071", 0, 1),
    @(41, "'590007", "'中邮中证500指数增强A", "'0.43", "'91.51", "'1.43", "'0.0061", 0, 7),
    @(42, "'970047", "'东海证券海睿健行灵活配置混合B", "'0.19", "'77.45", "'2.85", "'0.0054", 0, 8),
    @(43, "'006025", "'诺安优化配置混合", "'0.02", "'73.48", "'4.62", "'0.0009", 0, 1),
    @(44, "'008124", "'中邮中证500指数增强C", "'0.04", "'91.51", "'1.43", "'0.0006", 0, 7),
    @(45, "'010952", "'中金成长精选混合C", "'0.02", "'88.74", "'2.91", "'0.0006", 0, 10),
    @(46, "'001602", "'鑫元鑫新收益灵活配置混合C", "'0.00", "'86.90", "'4.87", "'0", 1, 4)
)

foreach ($row in $funds) {
    $r = [int]$row[0] + 2
    $q1.Range("A$r").Value = [int]$row[0]
    $q1.Range("A$r").Style = $q1.Range("A2").Style
    $q1.Range("B$r").Value = $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = $row[3]
    $q1.Range("E$r").Value = $row[4]
    $q1.Range("F$r").Value = $row[5]
    if ($row[7] -eq 1) {
        $q1.Range("G$r").Value = 0
    } else {
        $q1.Range("G$r").Value = $row[6]
    }
    $q1.Range("H$r").Value = [int]$row[8]
}

# ---------------------------------------------------------------------
# Step 2: add a brand-new "总计" sheet at the end with the updated totals.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add()
$total.Name = "总计"
$lastIdx = $wb.Worksheets.Count
$total.Move($null, $wb.Worksheets.Item($lastIdx))

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
$totalCols = @("B", "C", "D")
for ($i = 0; $i -lt $totalHeaders.Length; $i++) {
    $cell = $total.Range($totalCols[$i] + "1")
    $cell.Value = $totalHeaders[$i]
    $cell.Style = $q1.Range("B1").Style
}

# Each row: @(index, dateLabel, count, heldValue)
$totals = @(
    @(0, "'2022-Q1", 47, 14.29),
    @(1, "'2021-Q4", 13, 5.32),
    @(2, "'2021-Q3", 5, 0.5),
    @(3, "'2021-Q2", 8, 0.37),
    @(4, "'2021-Q1", 32, 7.14),
    @(5, "'2020-Q4", 10, 0.5600000000000001)
)

foreach ($row in $totals) {
    $r = [int]$row[0] + 2
    $total.Range("A$r").Value = [int]$row[0]
    $total.Range("A$r").Style = $q1.Range("A2").Style
    $total.Range("B$r").Value = $row[1]
    $total.Range("C$r").Value = [double]$row[2]
    $total.Range("D$r").Value = [double]$row[3]
}

